$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1794.1708
$ws.Range("I132").Value = 1905.0286
$ws.Range("K132").Value = 5715.085800000001
$ws.Range("M132").Value = -3185.085800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10334.8
$ws.Range("I32").Value = 6551.5977
$ws.Range("J32").Value = 27569.389
$ws.Range("K32").Value = 6551.5977
$ws.Range("L32").Value = 27569.389
$ws.Range("M32").Value = -6264.5977
$ws.Range("N32").Value = -28143.389
$ws.Range("H74").Value = 849
$ws.Range("I74").Value = 869.9286
$ws.Range("J74").Value = 769.0909
$ws.Range("K74").Value = 869.9286
$ws.Range("L74").Value = 769.0909
$ws.Range("M74").Value = 4.07140000000004
$ws.Range("N74").Value = -2517.0909
$ws.Range("H77").Value = 849
$ws.Range("I77").Value = 869.9286
$ws.Range("J77").Value = 769.0909
$ws.Range("K77").Value = 4349.643
$ws.Range("L77").Value = 3845.4545
$ws.Range("M77").Value = 18.35699999999997
$ws.Range("N77").Value = -12581.4545
$ws.Range("H88").Value = 2277.8948
$ws.Range("I88").Value = 1977.1428
$ws.Range("J88").Value = 3120
$ws.Range("K88").Value = 1977.1428
$ws.Range("L88").Value = 3120
$ws.Range("M88").Value = -1571.1428
$ws.Range("N88").Value = -3932
$ws.Range("H91").Value = 2277.8948
$ws.Range("I91").Value = 1977.1428
$ws.Range("J91").Value = 3120
$ws.Range("K91").Value = 1977.1428
$ws.Range("L91").Value = 3120
$ws.Range("M91").Value = -573.1428000000001
$ws.Range("N91").Value = -5928

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66669132
$ws.Range("I86").Value = 71430640
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 71430640
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -71429517
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 66669132
$ws.Range("I89").Value = 71430640
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 357153200
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -357147584
$ws.Range("N89").Value = -51232
$ws.Range("H94").Value = 925.2727
$ws.Range("I94").Value = 907.75
$ws.Range("J94").Value = 1100.5
$ws.Range("K94").Value = 907.75
$ws.Range("L94").Value = 1100.5
$ws.Range("M94").Value = -456.75
$ws.Range("N94").Value = -2002.5
$ws.Range("H134").Value = 70119.57000000001
$ws.Range("I134").Value = 4062.6667
$ws.Range("J134").Value = 169204.92
$ws.Range("K134").Value = 12188.0001
$ws.Range("L134").Value = 507614.76
$ws.Range("M134").Value = -9653.000100000001
$ws.Range("N134").Value = -512684.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2252663
$ws.Range("I31").Value = 2905750.8
$ws.Range("K31").Value = 2905750.8
$ws.Range("M31").Value = -2905455.8
$ws.Range("H34").Value = 2252663
$ws.Range("I34").Value = 2905750.8
$ws.Range("K34").Value = 2905750.8
$ws.Range("M34").Value = -2905548.8
$ws.Range("H94").Value = 3019.1282
$ws.Range("I94").Value = 3138.25
$ws.Range("J94").Value = 2988.3872
$ws.Range("K94").Value = 3138.25
$ws.Range("L94").Value = 2988.3872
$ws.Range("M94").Value = -2687.25
$ws.Range("N94").Value = -3890.3872
$ws.Range("H134").Value = 1673.1628
$ws.Range("I134").Value = 1634.6666
$ws.Range("J134").Value = 1762
$ws.Range("K134").Value = 4903.9998
$ws.Range("L134").Value = 5286
$ws.Range("M134").Value = -2368.9998
$ws.Range("N134").Value = -10356

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 485.68967
$ws.Range("I5").Value = 335.9091
$ws.Range("J5").Value = 956.4286
$ws.Range("K5").Value = 1007.7273
$ws.Range("L5").Value = 2869.2858
$ws.Range("M5").Value = -895.7273
$ws.Range("N5").Value = -3093.2858
$ws.Range("H122").Value = 526.5192
$ws.Range("I122").Value = 478
$ws.Range("J122").Value = 529.4897999999999
$ws.Range("K122").Value = 4302
$ws.Range("L122").Value = 4765.4082
$ws.Range("M122").Value = -1852
$ws.Range("N122").Value = -9665.4082
$ws.Range("H135").Value = 485.68967
$ws.Range("I135").Value = 335.9091
$ws.Range("J135").Value = 956.4286
$ws.Range("K135").Value = 3023.1819
$ws.Range("L135").Value = 8607.857399999999
$ws.Range("M135").Value = -488.1819
$ws.Range("N135").Value = -13677.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3315.2
$ws.Range("I80").Value = 3144
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 3144
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -2146
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 3315.2
$ws.Range("I83").Value = 3144
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 15720
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -10728
$ws.Range("N83").Value = -29984
$ws.Range("H102").Value = 1172.7059
$ws.Range("I102").Value = 1197.4
$ws.Range("J102").Value = 987.5
$ws.Range("K102").Value = 1197.4
$ws.Range("L102").Value = 987.5
$ws.Range("M102").Value = 424.5999999999999
$ws.Range("N102").Value = -4231.5
$ws.Range("H122").Value = 9274014
$ws.Range("I122").Value = 7983311
$ws.Range("J122").Value = 12500770
$ws.Range("K122").Value = 23949933
$ws.Range("L122").Value = 37502310
$ws.Range("M122").Value = -23947483
$ws.Range("N122").Value = -37507210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2776
$ws.Range("I68").Value = 2845
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 2845
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -2096
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 2776
$ws.Range("I71").Value = 2845
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 14225
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -10481
$ws.Range("N71").Value = -19988
$ws.Range("H132").Value = 1860.2712
$ws.Range("I132").Value = 1647.04
$ws.Range("J132").Value = 3044.889
$ws.Range("K132").Value = 4941.12
$ws.Range("L132").Value = 9134.667000000001
$ws.Range("M132").Value = -2411.12
$ws.Range("N132").Value = -14194.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1738.3684
$ws.Range("I81").Value = 1637.4166
$ws.Range("J81").Value = 1911.4286
$ws.Range("K81").Value = 3274.8332
$ws.Range("L81").Value = 3822.8572
$ws.Range("M81").Value = -2213.8332
$ws.Range("N81").Value = -5944.8572
$ws.Range("H84").Value = 1738.3684
$ws.Range("I84").Value = 1637.4166
$ws.Range("J84").Value = 1911.4286
$ws.Range("K84").Value = 16374.166
$ws.Range("L84").Value = 19114.286
$ws.Range("M84").Value = -11070.166
$ws.Range("N84").Value = -29722.286
$ws.Range("H110").Value = 28500
$ws.Range("J110").Value = 28500
$ws.Range("L110").Value = 28500
$ws.Range("N110").Value = -36680
$ws.Range("H132").Value = 1334.9395
$ws.Range("I132").Value = 1528.3077
$ws.Range("K132").Value = 4584.9231
$ws.Range("M132").Value = -2054.9231
$ws.Range("H136").Value = 523.1818
$ws.Range("I136").Value = 440.66666
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 1321.99998
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = 1228.00002
$ws.Range("N136").Value = -7200
